$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at 778, pushing existing rows 778-818 down to 779-819
$ws.Rows.Item(778).Insert()

# Populate the new row 778 with the data from the diff
$ws.Cells.Item(778, 1).Value = 3
$ws.Cells.Item(778, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(778, 3).Value = "Coquimbo"
$ws.Cells.Item(778, 4).Value = 45147
$ws.Cells.Item(778, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(778, 5).Value = 5
$ws.Cells.Item(778, 6).Value = 100112032
$ws.Cells.Item(778, 7).Value = "Zapallo italiano"
$ws.Cells.Item(778, 8).Value = "Sin especificar"
$ws.Cells.Item(778, 9).Value = "Primera"
$ws.Cells.Item(778, 10).Value = 115
$ws.Cells.Item(778, 11).Value = 13500
$ws.Cells.Item(778, 12).Value = 14000
$ws.Cells.Item(778, 13).Value = 13717
$ws.Cells.Item(778, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(778, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(778, 16).Value = 229
$ws.Cells.Item(778, 17).Value = 60
$ws.Cells.Item(778, 18).Value = "Hortaliza"
